$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell Z1: text date label matching style of N1:Y1 header cells
$z1 = $ws.Range("Z1")
$z1.Value = "'11-10-2020"
$z1.Font.Bold = $true
$z1.HorizontalAlignment = -4108
$z1.VerticalAlignment = -4160
$z1.Borders.LineStyle = 1

# Numeric data cells Z2:Z36
$ws.Range("Z2").Value = 3744
$ws.Range("Z3").Value = 697699
$ws.Range("Z4").Value = 9035
$ws.Range("Z5").Value = 163355
$ws.Range("Z6").Value = 183390
$ws.Range("Z7").Value = 11662
$ws.Range("Z8").Value = 111654
$ws.Range("Z9").Value = 3050
$ws.Range("Z10").Value = 278812
$ws.Range("Z11").Value = 32777
$ws.Range("Z12").Value = 130760
$ws.Range("Z13").Value = 128841
$ws.Range("Z14").Value = 14278
$ws.Range("Z15").Value = 70955
$ws.Range("Z16").Value = 82805
$ws.Range("Z17").Value = 569947
$ws.Range("Z18").Value = 182874
$ws.Range("Z19").Value = 3973
$ws.Range("Z20").Value = 127034
$ws.Range("Z21").Value = 1255779
$ws.Range("Z22").Value = 10396
$ws.Range("Z23").Value = 5045
$ws.Range("Z24").Value = 1984
$ws.Range("Z25").Value = 5694
$ws.Range("Z26").Value = 224273
$ws.Range("Z27").Value = 25955
$ws.Range("Z28").Value = 109767
$ws.Range("Z29").Value = 133918
$ws.Range("Z30").Value = 2816
$ws.Range("Z31").Value = 597033
$ws.Range("Z32").Value = 185128
$ws.Range("Z33").Value = 24086
$ws.Range("Z34").Value = 46470
$ws.Range("Z35").Value = 387149
$ws.Range("Z36").Value = 255838
